$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows to update: row number -> new EBITDA (column B) value. All dates in
# column A move from 2026/01/09 to 2026/01/10. Row 38 only has its date
# updated (EBITDA value stays the same).
$updates = @{
    2  = "7.96"
    8  = "8.46"
    14 = "3.14"
    20 = "13.49"
    26 = "11.64"
    32 = "28.91"
    38 = $null
    44 = "13.94"
    50 = "11.64"
    56 = "31.69"
    62 = "11.30"
    68 = "12.82"
    74 = "18.41"
}

foreach ($row in $updates.Keys) {
    # Column A holds the date as plain text (e.g. "2026/01/10"). Force the
    # cell to text format first so Excel does not auto-convert the string
    # into a date serial number, then restore the default "Normal" style
    # so no stray number formatting is left behind.
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026/01/10"
    $dateCell.Style = "Normal"

    $newB = $updates[$row]
    if ($newB -ne $null) {
        $ebitdaCell = $ws.Cells.Item($row, 2)
        $ebitdaCell.NumberFormat = "@"
        $ebitdaCell.Value = $newB
        $ebitdaCell.Style = "Normal"
    }
}
